# Update column F ("dSF") values to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F8").Value = -4
$ws.Range("F10").Value = -2
$ws.Range("F12").Value = 5
